$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 into the two new header cells (I1, J1)
# so they pick up the same bold/centered/bordered formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-34
$data = @(
    @(10, 11),
    @(9, 9),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(2, 2),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
